$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The post about "「忙しい」مشغول" (previously at row 565) was removed from
# the spreadsheet. Deleting the entire row shifts every subsequent row
# (566-694) up by one (now 565-693), matching the rest of the diff, and
# updates the sheet's used-range dimension from A1:C694 to A1:C693.
$ws.Rows.Item(565).Delete()
